# This workbook tracks weekly Kiwi price records. A new week of records
# (4 rows, one per "Calidad" grade) is inserted just above the existing
# block that starts at row 349, shifting all subsequent rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows above row 349 (keeps formatting of the row below).
$ws.Range("A349:A352").EntireRow.Insert()

# Common values shared by the 4 new rows.
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100101
$producto    = "Berries"
$categoriaId = 100101007
$categoria   = "Kiwi"
$variedad    = "Hayward"
$unidad      = "`$/bins (450 kilos)"
$kgUnidad    = 450
$fecha       = 44463

function Set-KiwiRow($row, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $origen, $precioKg) {
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-KiwiRow 349 "Especial"                28 330000 380000 350000 "Región Metropolitana" 778
Set-KiwiRow 350 "Extra (doble especial)"  28 430000 450000 434643 "Región Metropolitana" 966
Set-KiwiRow 351 "Primera"                 30 250000 270000 260000 "Región Metropolitana" 578
Set-KiwiRow 352 "Segunda"                 22 200000 200000 200000 "Región Metropolitana" 444
